$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QR")
$ws.Range("A1").Value = "TEST"
